$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("gens")
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("C17").Value = 0
$ws.Range("C25").Value = 400
$ws.Range("C26").Value = 50
$ws.Range("C28").Value = 50
$ws.Range("C29").Value = 50
$ws.Range("C30").Value = 50
$ws.Range("C31").Value = 50

$ws = $wb.Worksheets.Item("lines")
$ws.Range("C2").Value = 27.507133
$ws.Range("D2").Value = 0.028976112
$ws.Range("C3").Value = -96.96709300000001
$ws.Range("D3").Value = -0.90828046
$ws.Range("C4").Value = -38.54004
$ws.Range("D4").Value = 0.18996858
$ws.Range("C5").Value = -49.712129
$ws.Range("D5").Value = 0.083314972
$ws.Range("C6").Value = -19.780738
$ws.Range("D6").Value = 0.27143019
$ws.Range("C7").Value = 13.997097
$ws.Range("D7").Value = 1.088798
$ws.Range("C8").Value = -290.96419
$ws.Range("D8").Value = -1.1301536
$ws.Range("C9").Value = -123.71213
$ws.Range("D9").Value = 0.068226434
$ws.Range("C10").Value = -109.54004
$ws.Range("D10").Value = 0.19667335
$ws.Range("C11").Value = -155.78074
$ws.Range("D11").Value = 0.08623563200000001
$ws.Range("F12").Value = -9.1410573
$ws.Range("C13").Value = -24.565076
$ws.Range("D13").Value = -0.10306221
$ws.Range("C14").Value = 11.065076
$ws.Range("D14").Value = 0.10306221
$ws.Range("C15").Value = -149.94941
$ws.Range("D15").Value = 0.41759948
$ws.Range("C16").Value = -159.33069
$ws.Range("D16").Value = 0.35360174
$ws.Range("C17").Value = -219.93721
$ws.Range("D17").Value = 0.21147507
$ws.Range("C18").Value = -229.31849
$ws.Range("D18").Value = 0.14747733
$ws.Range("C19").Value = -219.0275
$ws.Range("D19").Value = 0.040522324
$ws.Range("C20").Value = -150.85913
$ws.Range("D20").Value = 0.27908024
$ws.Range("C21").Value = -202.61025
$ws.Range("D21").Value = 0.10452006
$ws.Range("C22").Value = -186.03893
$ws.Range("D22").Value = 0.36740939
$ws.Range("C23").Value = -95.63775
$ws.Range("D23").Value = 0.26288933
$ws.Range("C24").Value = -344.85913
$ws.Range("D24").Value = 0.39204129
$ws.Range("C25").Value = 115.19187
$ws.Range("D25").Value = 4.0072915
$ws.Range("C26").Value = -284.07803
$ws.Range("D26").Value = -6.1048423
$ws.Range("C27").Value = -284.07803
$ws.Range("D27").Value = -6.1048423
$ws.Range("C28").Value = 290.96419
$ws.Range("D28").Value = 0.6996189
$ws.Range("C29").Value = -450
$ws.Range("D29").Value = 6.4786082
$ws.Range("E29").Value = 22.410276
$ws.Range("C30").Value = 275.33274
$ws.Range("D30").Value = -0.15661717
$ws.Range("C31").Value = -307.93014
$ws.Range("D31").Value = 3.0175364
$ws.Range("C32").Value = -142.06986
$ws.Range("D32").Value = 3.5320874
$ws.Range("C33").Value = 5.1129567
$ws.Range("D33").Value = 2.801998
$ws.Range("C34").Value = 5.1129567
$ws.Range("D34").Value = 2.801998
$ws.Range("C35").Value = 47.16637
$ws.Range("D35").Value = -0.13618884
$ws.Range("C36").Value = 47.16637
$ws.Range("D36").Value = -0.13618884
$ws.Range("C37").Value = -16.83363
$ws.Range("D37").Value = -0.074903863
$ws.Range("C38").Value = -16.83363
$ws.Range("D38").Value = -0.074903863
$ws.Range("C39").Value = -157.93014
$ws.Range("D39").Value = -2.287447
$ws.Range("C40").Value = 344.65606
$ws.Range("D40").Value = [double]"-4.3485215e-17"

$ws = $wb.Worksheets.Item("bus")
$ws.Range("B2").Value = 108.85748
$ws.Range("C2").Value = 0.56874645
$ws.Range("B3").Value = 108.88645
$ws.Range("C3").Value = 0.56489545
$ws.Range("B4").Value = 107.9492
$ws.Range("C4").Value = 0.7733470099999999
$ws.Range("B5").Value = 108.96977
$ws.Range("C5").Value = 0.62802985
$ws.Range("B6").Value = 109.04745
$ws.Range("C6").Value = 0.60150548
$ws.Range("B7").Value = 109.15788
$ws.Range("C7").Value = 0.60287446
$ws.Range("C8").Value = 0.81223309
$ws.Range("B9").Value = 109.14106
$ws.Range("C9").Value = 0.71615809
$ws.Range("B10").Value = 109.038
$ws.Range("C10").Value = 0.75669047
$ws.Range("B11").Value = 109.24412
$ws.Range("C11").Value = 0.69790072
$ws.Range("B12").Value = 109.45559
$ws.Range("C12").Value = 0.8826479699999999
$ws.Range("B13").Value = 109.3916
$ws.Range("C13").Value = 0.89052825
$ws.Range("B14").Value = 109.49612
$ws.Range("C14").Value = 0.98778117
$ws.Range("B15").Value = 109.73467
$ws.Range("C15").Value = 0.94600881
$ws.Range("B16").Value = 106.11942
$ws.Range("C16").Value = 1.1690583
$ws.Range("B17").Value = 110.12672
$ws.Range("C17").Value = 1.1494757
$ws.Range("B18").Value = 94.195048
$ws.Range("C18").Value = 1.2664757
$ws.Range("B19").Value = 97.21258400000001
$ws.Range("C19").Value = 1.3095859
$ws.Range("B20").Value = 109.9701
$ws.Range("C20").Value = 1.0861492
$ws.Range("B21").Value = 109.83391
$ws.Range("C21").Value = 1.0672826
$ws.Range("B22").Value = 100.01458
$ws.Range("C22").Value = 1.3082565
$ws.Range("B23").Value = 97.727135
$ws.Range("C23").Value = 1.415649
$ws.Range("B24").Value = 109.75901
$ws.Range("C24").Value = 1.070986
$ws.Range("B25").Value = 106.81904
$ws.Range("C25").Value = 1.0177569
$ws.Range("B26").Value = 109.75901
$ws.Range("C26").Value = 1.039967
